$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether | Ether
$ws.Range("H15").Value = 1154
$ws.Range("I15").Value = 1154
$ws.Range("K15").Value = 3462
$ws.Range("M15").Value = -3293

# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 1065.0834
$ws.Range("I28").Value = 661.125
$ws.Range("J28").Value = 1873
$ws.Range("K28").Value = 661.125
$ws.Range("L28").Value = 1873
$ws.Range("M28").Value = -176.125
$ws.Range("N28").Value = -2843

# Row 53: No Accounting for Waste | Enchanted Electrum Ink
$ws.Range("H53").Value = 281.81818
$ws.Range("I53").Value = 281.81818
$ws.Range("K53").Value = 281.81818
$ws.Range("M53").Value = 355.18182

# Row 92: Whinier than the Sword | Enchanted Koppranickel Ink
$ws.Range("H92").Value = 365.92307
$ws.Range("I92").Value = 345.5
$ws.Range("J92").Value = 434
$ws.Range("K92").Value = 345.5
$ws.Range("L92").Value = 434
$ws.Range("M92").Value = 902.5
$ws.Range("N92").Value = -2930

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1124.5
$ws.Range("I98").Value = 599.4
$ws.Range("K98").Value = 599.4
$ws.Range("M98").Value = 898.6

# Row 107: Another Man's Ink | Enchanted Truegold Ink
$ws.Range("H107").Value = 281.81818
$ws.Range("I107").Value = 311.44446
$ws.Range("K107").Value = 311.44446
$ws.Range("M107").Value = 1608.55554

# Row 111: An Eye for Healing | Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 6749.25
$ws.Range("I111").Value = 6749.25
$ws.Range("K111").Value = 20247.75
$ws.Range("M111").Value = -17180.75

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 1177.8611
$ws.Range("J112").Value = 1547.9565
$ws.Range("L112").Value = 4643.8695
$ws.Range("N112").Value = -6859.8695

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1124.5
$ws.Range("I122").Value = 599.4
$ws.Range("K122").Value = 1798.2
$ws.Range("M122").Value = 651.8000000000002

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 2006.25
$ws.Range("I125").Value = 1998
$ws.Range("J125").Value = 2014.5
$ws.Range("K125").Value = 17982
$ws.Range("L125").Value = 18130.5
$ws.Range("M125").Value = -15522
$ws.Range("N125").Value = -23050.5

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 2119.125
$ws.Range("I135").Value = 2119.125
$ws.Range("K135").Value = 19072.125
$ws.Range("M135").Value = -16537.125

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2937.8823
$ws.Range("I137").Value = 745.625
$ws.Range("J137").Value = 4886.5557
$ws.Range("K137").Value = 2236.875
$ws.Range("L137").Value = 14659.6671
$ws.Range("M137").Value = 313.125
$ws.Range("N137").Value = -19759.6671

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2675.8
$ws.Range("J138").Value = 2958.5217
$ws.Range("L138").Value = 8875.5651
$ws.Range("N138").Value = -19155.5651

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 9153.6875
$ws.Range("I32").Value = 9153.6875
$ws.Range("K32").Value = 9153.6875
$ws.Range("M32").Value = -8866.6875

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1388.8
$ws.Range("I45").Value = 1388.8
$ws.Range("K45").Value = 1388.8
$ws.Range("M45").Value = -1011.8

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2180
$ws.Range("I61").Value = 1580.875
$ws.Range("J61").Value = 3138.6
$ws.Range("K61").Value = 1580.875
$ws.Range("L61").Value = 3138.6
$ws.Range("M61").Value = -1368.875
$ws.Range("N61").Value = -3562.6

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 2001
$ws.Range("I74").Value = 2001
$ws.Range("K74").Value = 2001
$ws.Range("M74").Value = -1127

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 2001
$ws.Range("I77").Value = 2001
$ws.Range("K77").Value = 10005
$ws.Range("M77").Value = -5637

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1100
$ws.Range("I110").Value = 1100
$ws.Range("K110").Value = 1100
$ws.Range("M110").Value = 945

# Row 131: Additions to the Armoire | Chondrite Top of Maiming
$ws.Range("H131").Value = 54949
$ws.Range("J131").Value = 54949
$ws.Range("L131").Value = 54949
$ws.Range("N131").Value = -65029

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2180
$ws.Range("I136").Value = 1580.875
$ws.Range("J136").Value = 3138.6
$ws.Range("K136").Value = 4742.625
$ws.Range("L136").Value = 9415.799999999999
$ws.Range("M136").Value = -2192.625
$ws.Range("N136").Value = -14515.8

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 3466.6667
$ws.Range("I105").Value = 1401
$ws.Range("J105").Value = 4499.5
$ws.Range("K105").Value = 1401
$ws.Range("L105").Value = 4499.5
$ws.Range("M105").Value = 346
$ws.Range("N105").Value = -7993.5

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 1490.25
$ws.Range("I107").Value = 1402
$ws.Range("J107").Value = 2108
$ws.Range("K107").Value = 1402
$ws.Range("L107").Value = 2108
$ws.Range("M107").Value = 518
$ws.Range("N107").Value = -5948

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1788.4286
$ws.Range("I134").Value = 1591.6471
$ws.Range("K134").Value = 4774.9413
$ws.Range("M134").Value = -2239.9413

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 2184
$ws.Range("I16").Value = 2184
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2184
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1897
$ws.Range("N16").ClearContents()

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 3828.1538
$ws.Range("I31").Value = 2225.4
$ws.Range("J31").Value = 9170.666999999999
$ws.Range("K31").Value = 2225.4
$ws.Range("L31").Value = 9170.666999999999
$ws.Range("M31").Value = -1930.4
$ws.Range("N31").Value = -9760.666999999999

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 3828.1538
$ws.Range("I34").Value = 2225.4
$ws.Range("J34").Value = 9170.666999999999
$ws.Range("K34").Value = 2225.4
$ws.Range("L34").Value = 9170.666999999999
$ws.Range("M34").Value = -2023.4
$ws.Range("N34").Value = -9574.666999999999

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2580.25
$ws.Range("I58").Value = 2138.7144
$ws.Range("J58").Value = 2923.6667
$ws.Range("K58").Value = 2138.7144
$ws.Range("L58").Value = 2923.6667
$ws.Range("M58").Value = -1935.7144
$ws.Range("N58").Value = -3329.6667

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 1097.6154
$ws.Range("I107").Value = 553.3333
$ws.Range("K107").Value = 553.3333
$ws.Range("M107").Value = 1366.6667

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 2184
$ws.Range("I113").Value = 2184
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2184
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -14
$ws.Range("N113").ClearContents()

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 3310.7778
$ws.Range("I132").Value = 2599.6667
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 7799.000100000001
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -5269.000100000001
$ws.Range("N132").Value = -16058.9999

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2580.25
$ws.Range("I136").Value = 2138.7144
$ws.Range("J136").Value = 2923.6667
$ws.Range("K136").Value = 6416.1432
$ws.Range("L136").Value = 8771.000100000001
$ws.Range("M136").Value = -3866.1432
$ws.Range("N136").Value = -13871.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 72689
$ws.Range("I122").Value = 1030.1428
$ws.Range("J122").Value = 144347.86
$ws.Range("K122").Value = 9271.2852
$ws.Range("L122").Value = 1299130.74
$ws.Range("M122").Value = -6821.2852
$ws.Range("N122").Value = -1304030.74

# Row 129: Comfort Food | Yakow Moussaka
$ws.Range("H129").Value = 1436.625
$ws.Range("J129").Value = 1999.3334
$ws.Range("L129").Value = 5998.0002
$ws.Range("N129").Value = -15998.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 1186.7142
$ws.Range("I107").Value = 1186.7142
$ws.Range("K107").Value = 1186.7142
$ws.Range("M107").Value = 733.2858000000001

# Row 129: The Needle That Binds | Manganese Needle
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2200.8948
$ws.Range("I132").Value = 1116.25
$ws.Range("K132").Value = 3348.75
$ws.Range("M132").Value = -818.75

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 2986.5
$ws.Range("I132").Value = 2090.6667
$ws.Range("K132").Value = 6272.000100000001
$ws.Range("M132").Value = -3742.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 1051.5
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1103
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 3309
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -7149

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 974.6429000000001
$ws.Range("I113").Value = 1004.25
$ws.Range("K113").Value = 3012.75
$ws.Range("M113").Value = -842.75

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 2052.6667
$ws.Range("I126").Value = 2052.6667
$ws.Range("K126").Value = 6158.000100000001
$ws.Range("M126").Value = -3688.000100000001

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 3476.238
$ws.Range("I132").Value = 3256.2222
$ws.Range("J132").Value = 3641.25
$ws.Range("K132").Value = 9768.6666
$ws.Range("L132").Value = 10923.75
$ws.Range("M132").Value = -7238.6666
$ws.Range("N132").Value = -15983.75
